$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new activity row (row 4) with the latest entry.
$ws.Range("B4").Value = "21/06/2016"
$ws.Range("C4").Value = "Realización de documento funciones asesores"
$ws.Range("D4").Value = "Análisis, Requisitos"
$ws.Range("E4").Value = "1. Se crea documento de word con las funciones resultantes de las encuestas a los asesores.`n2. se crea una carpeta de requisitos, donde se iran guardando los documentos de requisitos.`n"

# Match row height used for the other detailed entry.
$ws.Range("B4:F4").RowHeight = 75

# Move the active selection down to the next empty row, mirroring the author's workflow.
$ws.Range("B5").Select() | Out-Null
